$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in C1 and D1
$ws.Range("C1").Value = "Taken By"
$ws.Range("D1").Value = "Received"

# Row 2: verifone / verifone
$ws.Range("C2").Value = "רכעקכעקכ"
$ws.Range("D2").Value = "Yes"

# Row 3: aarar / moataz - clear the "taken by" note, keep/move "Yes" to D3
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = "Yes"

# Row 4: abdush / shiran
$ws.Range("C4").Value = "DDDD"
$ws.Range("D4").Value = "Yes"

# Row 5: abu shamala / zahiya
$ws.Range("C5").Value = "DDDDKJYKJHJ"
$ws.Range("D5").Value = "Yes"

# Row 7: adani / inbal
$ws.Range("C7").Value = "DVFGHDFVSDGNFGHGH"
$ws.Range("D7").Value = "Yes"

# Row 57: dahan / limor - clear "Yes" from C57
$ws.Range("C57").Value = ""

# Row 96: gonsharovich / eldad - clear both
$ws.Range("C96").Value = ""
$ws.Range("D96").Value = ""

# Row 118: irmiyahoo / liat - clear both
$ws.Range("C118").Value = ""
$ws.Range("D118").Value = ""

# Row 254: shmartov / nastya - clear both
$ws.Range("C254").Value = ""
$ws.Range("D254").Value = ""

# Row 261: sinenco / peter - clear both
$ws.Range("C261").Value = ""
$ws.Range("D261").Value = ""
